# Workbook/worksheet references (workbook already open as ActiveWorkbook)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# The "Recorded By" column (G) contains entries formatted as "Administrator, <name>".
# Re-order every such entry to "<name>, Administrator".
$oldText = "Administrator, Miss Dina Nasr"
$newText = "Miss Dina Nasr, Administrator"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
